{"js": "// Update the date heading and the 100 arithmetic-problem cells in the\n// table, in document order. The source doc has one duplicate problem\n// text (\"16+61=\") that maps to two different replacements depending on\n// position, so replacement must be positional (not a global find/replace).\n\nconst newCellValues = [\"48-15=\", \"14+16=\", \"7+45=\", \"96-70=\", \"61+11=\", \"69-4=\", \"1+79=\", \"96-29=\", \"48-21=\", \"10-5=\", \"78-63=\", \"16+18=\", \"30+34=\", \"31+7=\", \"22+42=\", \"69-2=\", \"78+6=\", \"27+57=\", \"74+10=\", \"10+23=\", \"45+46=\", \"7+16=\", \"83+8=\", \"47+3=\", \"48+12=\", \"28+10=\", \"96-88=\", \"87-85=\", \"72+16=\", \"52-33=\", \"23+9=\", \"29-26=\", \"35+55=\", \"61+14=\", \"37+26=\", \"35+30=\", \"85-81=\", \"76-42=\", \"49+27=\", \"98-85=\", \"28+27=\", \"20+15=\", \"65+16=\", \"63-6=\", \"21+60=\", \"90-31=\", \"17+41=\", \"35+23=\", \"36+51=\", \"47+29=\", \"51-12=\", \"29+14=\", \"50+38=\", \"78+15=\", \"14+2=\", \"66-48=\", \"98-55=\", \"43-6=\", \"36+0=\", \"22+33=\", \"42+48=\", \"98-12=\", \"37+11=\", \"9+10=\", \"22+69=\", \"68-50=\", \"52+15=\", \"64-55=\", \"53-10=\", \"8+62=\", \"45+52=\", \"35-16=\", \"21+71=\", \"51+27=\", \"57+26=\", \"17+59=\", \"98-90=\", \"33+25=\", \"45-42=\", \"63+25=\", \"92-39=\", \"89-53=\", \"23+69=\", \"96-5=\", \"24+21=\", \"77-25=\", \"52-48=\", \"36+62=\", \"1+71=\", \"15+81=\", \"71-51=\", \"61+13=\", \"81-70=\", \"6+50=\", \"6+56=\", \"86-20=\", \"43+8=\", \"59-26=\", \"3+94=\", \"55-43=\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// First paragraph of the document holds the date heading.\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(\"2023-05-18 Thursday\", \"Replace\");\n\n// The single table holds the 20x5 grid of arithmetic problems.\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet idx = 0;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (const cell of cells.items) {\n    cell.value = newCellValues[idx];\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 100 arithmetic-problem cells in the\n# table, in document order. The source doc has one duplicate problem\n# text (\"16+61=\") that maps to two different replacements depending on\n# position, so replacement must be positional (Cell(r,c)), not a global\n# Find/Replace.\n\n$d = $word.ActiveDocument\n\n# First paragraph of the document holds the date heading.\n$d.Paragraphs(1).Range.Text = \"2023-05-18 Thursday\"\n\n# New values for the 20x5 grid, in row-major (top-to-bottom, left-to-right) order.\n$newValues = @(\n    @(\"48-15=\", \"14+16=\", \"7+45=\", \"96-70=\", \"61+11=\"),\n    @(\"69-4=\", \"1+79=\", \"96-29=\", \"48-21=\", \"10-5=\"),\n    @(\"78-63=\", \"16+18=\", \"30+34=\", \"31+7=\", \"22+42=\"),\n    @(\"69-2=\", \"78+6=\", \"27+57=\", \"74+10=\", \"10+23=\"),\n    @(\"45+46=\", \"7+16=\", \"83+8=\", \"47+3=\", \"48+12=\"),\n    @(\"28+10=\", \"96-88=\", \"87-85=\", \"72+16=\", \"52-33=\"),\n    @(\"23+9=\", \"29-26=\", \"35+55=\", \"61+14=\", \"37+26=\"),\n    @(\"35+30=\", \"85-81=\", \"76-42=\", \"49+27=\", \"98-85=\"),\n    @(\"28+27=\", \"20+15=\", \"65+16=\", \"63-6=\", \"21+60=\"),\n    @(\"90-31=\", \"17+41=\", \"35+23=\", \"36+51=\", \"47+29=\"),\n    @(\"51-12=\", \"29+14=\", \"50+38=\", \"78+15=\", \"14+2=\"),\n    @(\"66-48=\", \"98-55=\", \"43-6=\", \"36+0=\", \"22+33=\"),\n    @(\"42+48=\", \"98-12=\", \"37+11=\", \"9+10=\", \"22+69=\"),\n    @(\"68-50=\", \"52+15=\", \"64-55=\", \"53-10=\", \"8+62=\"),\n    @(\"45+52=\", \"35-16=\", \"21+71=\", \"51+27=\", \"57+26=\"),\n    @(\"17+59=\", \"98-90=\", \"33+25=\", \"45-42=\", \"63+25=\"),\n    @(\"92-39=\", \"89-53=\", \"23+69=\", \"96-5=\", \"24+21=\"),\n    @(\"77-25=\", \"52-48=\", \"36+62=\", \"1+71=\", \"15+81=\"),\n    @(\"71-51=\", \"61+13=\", \"81-70=\", \"6+50=\", \"6+56=\"),\n    @(\"86-20=\", \"43+8=\", \"59-26=\", \"3+94=\", \"55-43=\")\n)\n\n$t = $d.Tables(1)\nfor ($r = 1; $r -le 20; $r++) {\n    for ($c = 1; $c -le 5; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
